# Auto-generated edit script: updates khl referees stats (Главные & Линейные sheets)
# Refreshes per-row stat columns (C..K, plus a few V-column tweaks) and the
# as_of_utc timestamp in column AA for every data row, matching the upstream diff.

$wb = $excel.ActiveWorkbook

$wsMain   = $wb.Worksheets.Item(2)   # "Главные"
$wsLinear = $wb.Worksheets.Item(3)   # "Линейные"


# ---- Главные ----
$wsMain.Range("AA2").Value = "2025-11-28 03:05:41"
$wsMain.Range("C3").Value = 29
$wsMain.Range("D3").Value = 519
$wsMain.Range("E3").Value = 233
$wsMain.Range("F3").Value = 286
$wsMain.Range("G3").Value = 17.9
$wsMain.Range("H3").Value = 8.029999999999999
$wsMain.Range("I3").Value = 9.859999999999999
$wsMain.Range("J3").Value = 114
$wsMain.Range("K3").Value = 118
$wsMain.Range("AA3").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA4").Value = "2025-11-28 03:05:41"
$wsMain.Range("C5").Value = 29
$wsMain.Range("D5").Value = 471
$wsMain.Range("E5").Value = 240
$wsMain.Range("F5").Value = 231
$wsMain.Range("G5").Value = 16.24
$wsMain.Range("H5").Value = 8.279999999999999
$wsMain.Range("I5").Value = 7.97
$wsMain.Range("J5").Value = 115
$wsMain.Range("K5").Value = 108
$wsMain.Range("AA5").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA6").Value = "2025-11-28 03:05:41"
$wsMain.Range("C7").Value = 18
$wsMain.Range("D7").Value = 239
$wsMain.Range("E7").Value = 104
$wsMain.Range("F7").Value = 135
$wsMain.Range("G7").Value = 13.28
$wsMain.Range("H7").Value = 5.78
$wsMain.Range("I7").Value = 7.5
$wsMain.Range("J7").Value = 52
$wsMain.Range("K7").Value = 50
$wsMain.Range("AA7").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA8").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA9").Value = "2025-11-28 03:05:41"
$wsMain.Range("C10").Value = 19
$wsMain.Range("D10").Value = 326
$wsMain.Range("E10").Value = 164
$wsMain.Range("F10").Value = 162
$wsMain.Range("G10").Value = 17.16
$wsMain.Range("H10").Value = 8.630000000000001
$wsMain.Range("I10").Value = 8.529999999999999
$wsMain.Range("J10").Value = 82
$wsMain.Range("K10").Value = 71
$wsMain.Range("AA10").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA11").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA12").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA13").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA14").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA15").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA16").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA17").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA18").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA19").Value = "2025-11-28 03:05:41"
$wsMain.Range("C20").Value = 28
$wsMain.Range("D20").Value = 487
$wsMain.Range("E20").Value = 213
$wsMain.Range("F20").Value = 274
$wsMain.Range("G20").Value = 17.39
$wsMain.Range("H20").Value = 7.61
$wsMain.Range("I20").Value = 9.789999999999999
$wsMain.Range("J20").Value = 99
$wsMain.Range("K20").Value = 102
$wsMain.Range("V20").Value = 14
$wsMain.Range("AA20").Value = "2025-11-28 03:05:41"
$wsMain.Range("C21").Value = 25
$wsMain.Range("D21").Value = 346
$wsMain.Range("E21").Value = 150
$wsMain.Range("F21").Value = 196
$wsMain.Range("G21").Value = 13.84
$wsMain.Range("H21").Value = 6
$wsMain.Range("I21").Value = 7.84
$wsMain.Range("J21").Value = 65
$wsMain.Range("K21").Value = 83
$wsMain.Range("AA21").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA22").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA23").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA24").Value = "2025-11-28 03:05:41"
$wsMain.Range("C25").Value = 29
$wsMain.Range("D25").Value = 463
$wsMain.Range("E25").Value = 226
$wsMain.Range("F25").Value = 237
$wsMain.Range("G25").Value = 15.97
$wsMain.Range("I25").Value = 8.17
$wsMain.Range("J25").Value = 108
$wsMain.Range("K25").Value = 111
$wsMain.Range("V25").Value = 10
$wsMain.Range("AA25").Value = "2025-11-28 03:05:41"
$wsMain.Range("AA26").Value = "2025-11-28 03:05:41"

# ---- Линейные ----
$wsLinear.Range("AA2").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA3").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA4").Value = "2025-11-28 03:05:41"
$wsLinear.Range("C5").Value = 15
$wsLinear.Range("D5").Value = 214
$wsLinear.Range("E5").Value = 114
$wsLinear.Range("F5").Value = 100
$wsLinear.Range("G5").Value = 14.27
$wsLinear.Range("H5").Value = 7.6
$wsLinear.Range("I5").Value = 6.67
$wsLinear.Range("J5").Value = 57
$wsLinear.Range("K5").Value = 50
$wsLinear.Range("AA5").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA6").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA7").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA8").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA9").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA10").Value = "2025-11-28 03:05:41"
$wsLinear.Range("C11").Value = 19
$wsLinear.Range("D11").Value = 270
$wsLinear.Range("E11").Value = 118
$wsLinear.Range("F11").Value = 152
$wsLinear.Range("G11").Value = 14.21
$wsLinear.Range("H11").Value = 6.21
$wsLinear.Range("I11").Value = 8
$wsLinear.Range("J11").Value = 59
$wsLinear.Range("K11").Value = 71
$wsLinear.Range("AA11").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA12").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA13").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA14").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA15").Value = "2025-11-28 03:05:41"
$wsLinear.Range("C16").Value = 27
$wsLinear.Range("D16").Value = 471
$wsLinear.Range("E16").Value = 218
$wsLinear.Range("F16").Value = 253
$wsLinear.Range("G16").Value = 17.44
$wsLinear.Range("H16").Value = 8.07
$wsLinear.Range("I16").Value = 9.369999999999999
$wsLinear.Range("J16").Value = 99
$wsLinear.Range("K16").Value = 109
$wsLinear.Range("AA16").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA17").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA18").Value = "2025-11-28 03:05:41"
$wsLinear.Range("C19").Value = 27
$wsLinear.Range("D19").Value = 447
$wsLinear.Range("E19").Value = 216
$wsLinear.Range("F19").Value = 231
$wsLinear.Range("G19").Value = 16.56
$wsLinear.Range("I19").Value = 8.56
$wsLinear.Range("J19").Value = 103
$wsLinear.Range("K19").Value = 103
$wsLinear.Range("V19").Value = 12
$wsLinear.Range("AA19").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA20").Value = "2025-11-28 03:05:41"
$wsLinear.Range("C21").Value = 31
$wsLinear.Range("D21").Value = 603
$wsLinear.Range("E21").Value = 250
$wsLinear.Range("F21").Value = 353
$wsLinear.Range("G21").Value = 19.45
$wsLinear.Range("H21").Value = 8.06
$wsLinear.Range("I21").Value = 11.39
$wsLinear.Range("J21").Value = 115
$wsLinear.Range("K21").Value = 144
$wsLinear.Range("AA21").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA22").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA23").Value = "2025-11-28 03:05:41"
$wsLinear.Range("C24").Value = 30
$wsLinear.Range("D24").Value = 534
$wsLinear.Range("E24").Value = 211
$wsLinear.Range("F24").Value = 323
$wsLinear.Range("G24").Value = 17.8
$wsLinear.Range("H24").Value = 7.03
$wsLinear.Range("I24").Value = 10.77
$wsLinear.Range("J24").Value = 93
$wsLinear.Range("K24").Value = 124
$wsLinear.Range("AA24").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA25").Value = "2025-11-28 03:05:41"
$wsLinear.Range("AA26").Value = "2025-11-28 03:05:41"
